$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "69.152.31"
Set-TextValue 2 5 "  +0.72%  "

# Row 3
Set-TextValue 3 4 "3.769.33"
Set-TextValue 3 5 "  +0.45%  "

# Row 4
Set-TextValue 4 5 "  -0.08%  "

# Row 5
Set-TextValue 5 4 "602.46"
Set-TextValue 5 5 "  +0.20%  "

# Row 6
Set-TextValue 6 4 "166.26"
Set-TextValue 6 5 "  -1.41%  "

# Row 7
Set-TextValue 7 4 "3.763.92"
Set-TextValue 7 5 "  +0.36%  "

# Row 8
Set-TextValue 8 5 "  -0.01%  "

# Row 9
Set-TextValue 9 5 "  +0.56%  "

# Row 10
Set-TextValue 10 5 "  +4.72%  "

# Row 11
Set-TextValue 11 5 "  +0.72%  "

# Row 12
Set-TextValue 12 4 "0.460"
Set-TextValue 12 5 "  -0.40%  "

# Row 13
Set-TextValue 13 4 "37.72"
Set-TextValue 13 5 "  -1.32%  "

# Row 14
Set-TextValue 14 5 "  +0.56%  "

# Row 15
Set-TextValue 15 4 "4.392.75"
Set-TextValue 15 5 "  +0.35%  "

# Row 16
Set-TextValue 16 4 "3.789.11"
Set-TextValue 16 5 "  +0.89%  "

# Row 17
Set-TextValue 17 4 "69.238.17"
Set-TextValue 17 5 "  +0.80%  "

# Row 18
Set-TextValue 18 4 "7.41"
Set-TextValue 18 5 "  +1.69%  "

# Row 19
Set-TextValue 19 4 "17.65"
Set-TextValue 19 5 "  +3.21%  "

# Row 20
Set-TextValue 20 5 "  -1.03%  "

# Row 21
Set-TextValue 21 4 "11.18"
Set-TextValue 21 5 "  +2.54%  "

# Row 22
Set-TextValue 22 4 "493.49"
Set-TextValue 22 5 "  -0.32%  "

# Row 23
Set-TextValue 23 4 "0.725"
Set-TextValue 23 5 "  -0.38%  "

# Row 24
Set-TextValue 24 4 "0.0000151"
Set-TextValue 24 5 "  -1.73%  "

# Row 25
Set-TextValue 25 4 "84.80"
Set-TextValue 25 5 "  -0.53%  "

# Row 27
Set-TextValue 27 4 "12.28"
Set-TextValue 27 5 "  -0.62%  "

# Row 28
Set-TextValue 28 2 "RenderToken"
Set-TextValue 28 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 28 4 "10.08"
Set-TextValue 28 5 "  -1.53%  "

# Row 29
Set-TextValue 29 2 "Dai"
Set-TextValue 29 3 "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue 29 4 "1.00"
Set-TextValue 29 5 "  +0.00%  "

# Row 30
Set-TextValue 30 2 "PancakeSwap"
Set-TextValue 30 3 "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue 30 4 "2.98"
Set-TextValue 30 5 "  -0.25%  "

# Row 31
Set-TextValue 31 2 "NEARProtocol"
Set-TextValue 31 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue 31 4 "8.16"
Set-TextValue 31 5 "  +3.24%  "

# Row 32
Set-TextValue 32 2 "ImmutableX"
Set-TextValue 32 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue 32 4 "2.43"
Set-TextValue 32 5 "  -3.90%  "

# Row 33
Set-TextValue 33 2 "EthereumClassic"
Set-TextValue 33 3 "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue 33 4 "31.88"
Set-TextValue 33 5 "  +0.02%  "

# Row 34
Set-TextValue 34 2 "WrappedeETH"
Set-TextValue 34 3 "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue 34 4 "3.910.97"
Set-TextValue 34 5 "  +0.40%  "

# Row 35
Set-TextValue 35 2 "RenzoRestakedETH"
Set-TextValue 35 3 "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue 35 4 "3.720.01"
Set-TextValue 35 5 "  +0.95%  "

# Row 36
Set-TextValue 36 2 "Hedera"
Set-TextValue 36 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue 36 4 "0.108"
Set-TextValue 36 5 "  -0.65%  "

# Row 37
Set-TextValue 37 2 "Filecoin"
Set-TextValue 37 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue 37 4 "5.96"
Set-TextValue 37 5 "  +1.78%  "

# Row 38
Set-TextValue 38 2 "Mantle"
Set-TextValue 38 3 "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue 38 4 "1.01"
Set-TextValue 38 5 "  -0.52%  "

# Row 39
Set-TextValue 39 2 "Kaspa"
Set-TextValue 39 3 "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue 39 4 "0.138"
Set-TextValue 39 5 "  +3.93%  "

# Row 40
Set-TextValue 40 2 "FirstDigitalUSD"
Set-TextValue 40 3 "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue 40 4 "1.00"
Set-TextValue 40 5 "  +0.01%  "

# Row 41
Set-TextValue 41 2 "dogwifhat"
Set-TextValue 41 3 "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue 41 4 "3.10"
Set-TextValue 41 5 "  +6.41%  "

# Row 42
Set-TextValue 42 2 "TheGraph"
Set-TextValue 42 3 "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue 42 4 "0.325"
Set-TextValue 42 5 "  +0.39%  "

# Row 43
Set-TextValue 43 2 "Stacks"
Set-TextValue 43 3 "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue 43 4 "1.99"
Set-TextValue 43 5 "  +1.01%  "

# Row 44
Set-TextValue 44 2 "OKB"
Set-TextValue 44 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue 44 4 "48.52"
Set-TextValue 44 5 "  -0.55%  "

# Row 45
Set-TextValue 45 2 "Bittensor"
Set-TextValue 45 3 "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue 45 4 "426.03"
Set-TextValue 45 5 "  -3.24%  "

# Row 46
Set-TextValue 46 2 "Cosmos"
Set-TextValue 46 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue 46 4 "8.44"
Set-TextValue 46 5 "  -0.38%  "

# Row 47
Set-TextValue 47 2 "USDe"
Set-TextValue 47 3 "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue 47 4 "1.00"
Set-TextValue 47 5 "  +0.02%  "

# Row 48
Set-TextValue 48 2 "Arweave"
Set-TextValue 48 3 "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue 48 4 "40.05"
Set-TextValue 48 5 "  -0.68%  "

# Row 49
Set-TextValue 49 2 "Monero"
Set-TextValue 49 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue 49 4 "142.20"
Set-TextValue 49 5 "  +0.40%  "

# Row 50
Set-TextValue 50 2 "Maker"
Set-TextValue 50 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue 50 4 "2.811.33"
Set-TextValue 50 5 "  +0.53%  "

# Row 51
Set-TextValue 51 2 "ONDO"
Set-TextValue 51 3 "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue 51 4 "1.29"
Set-TextValue 51 5 "  +8.48%  "
